$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header summary numbers (Valor Mora total, Cant. Trabajadores, Cant. Periodos) ---
$ws.Range("E11").Value = 907902
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 6

# --- Refresh worker detail table (rows 16-24) with the new period-of-account data ---
# Row 16: CARLOS JAVIER JARABA GUTIERREZ - periodo 2502
$ws.Range("C16").Value = "1052984539"
$ws.Range("D16").Value = "CARLOS JAVIER JARABA GUTIERREZ"
$ws.Range("E16").Value = "2502"
$ws.Range("F16").Value = 26572
$ws.Range("G16").Value = 1423500

# Row 17: CARLOS JAVIER JARABA GUTIERREZ - periodo 2503
$ws.Range("C17").Value = "1052984539"
$ws.Range("D17").Value = "CARLOS JAVIER JARABA GUTIERREZ"
$ws.Range("E17").Value = "2503"
$ws.Range("F17").Value = 49348
$ws.Range("G17").Value = 1423500

# Row 18: CARLOS JAVIER JARABA GUTIERREZ - periodo 2503
$ws.Range("C18").Value = "1052984539"
$ws.Range("D18").Value = "CARLOS JAVIER JARABA GUTIERREZ"
$ws.Range("E18").Value = "2503"
$ws.Range("F18").Value = 32266
$ws.Range("G18").Value = 1423500

# Row 19: IVAN DAVID ATENCIA DELGADO - periodo 2505
$ws.Range("C19").Value = "1002497743"
$ws.Range("D19").Value = "IVAN DAVID ATENCIA DELGADO"
$ws.Range("E19").Value = "2505"
$ws.Range("F19").Value = 22776
$ws.Range("G19").Value = 1423500

# Row 20: JOSE FERNANDO GAITAN GAITAN - periodo 2505 (unchanged values, kept for completeness)
$ws.Range("C20").Value = "1116043174"
$ws.Range("D20").Value = "JOSE FERNANDO GAITAN GAITAN"
$ws.Range("E20").Value = "2505"
$ws.Range("F20").Value = 180000
$ws.Range("G20").Value = 4500000

# Row 21: IVAN DAVID ATENCIA DELGADO - periodo 2506
$ws.Range("C21").Value = "1002497743"
$ws.Range("D21").Value = "IVAN DAVID ATENCIA DELGADO"
$ws.Range("E21").Value = "2506"
$ws.Range("F21").Value = 56940
$ws.Range("G21").Value = 1423500

# Row 22: JOSE FERNANDO GAITAN GAITAN - periodo 2506
$ws.Range("C22").Value = "1116043174"
$ws.Range("D22").Value = "JOSE FERNANDO GAITAN GAITAN"
$ws.Range("E22").Value = "2506"
$ws.Range("F22").Value = 180000
$ws.Range("G22").Value = 4500000

# Row 23: JOSE FERNANDO GAITAN GAITAN - periodo 2507
$ws.Range("C23").Value = "1116043174"
$ws.Range("D23").Value = "JOSE FERNANDO GAITAN GAITAN"
$ws.Range("E23").Value = "2507"
$ws.Range("F23").Value = 180000
$ws.Range("G23").Value = 4500000

# Row 24: JOSE FERNANDO GAITAN GAITAN - periodo 2508 (replaces the old KEVIN ALBERTO ECHEVERRIA row)
$ws.Range("C24").Value = "1116043174"
$ws.Range("D24").Value = "JOSE FERNANDO GAITAN GAITAN"
$ws.Range("E24").Value = "2508"
$ws.Range("F24").Value = 180000
$ws.Range("G24").Value = 4500000
